# Add "Clinical Trial, n (%)" breakdown rows (20-29) for the 1399 Dx cohort
# (and partial data for the Px / Validation cohorts), per commit:
# "add trial info for 1399 Dx cohort."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge the new category label cell (A20:A29) BEFORE copying styles onto the
# range - merging first keeps a uniform border style instead of Excel
# splitting the border into top/middle/bottom pieces across the merged cells.
$ws.Range("A20:A29").Merge()

# Copy the existing bold/centered/bordered look from the last category block
# (A18:B19 -> Treatment Arm row pair) across the whole new block so that the
# new cells pick up the same style index as the rest of the table.
$ws.Range("A18:B19").Copy()
$ws.Range("A20:B29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row/category label
$ws.Range("A20").Value2 = "Clinical Trial, n (%)"

# Trial names (column B)
$ws.Range("B20").Value2 = "AAML03P1"
$ws.Range("B21").Value2 = "AAML0531"
$ws.Range("B22").Value2 = "AAML1031"
$ws.Range("B23").Value2 = "Beat AML Consortium"
$ws.Range("B24").Value2 = "CCG2961"
$ws.Range("B25").Value2 = "CETLAM SMD-09 (MDS-tAML)"
$ws.Range("B26").Value2 = "Japanese AML05"
$ws.Range("B27").Value2 = "TCGA AML"
$ws.Range("B28").Value2 = "AML02"
$ws.Range("B29").Value2 = "AML08"

# MethylScoreAML Dx Discovery cohort (column C) - 1399 n
$ws.Range("C20").Value2 = "60 (4.3)"
$ws.Range("C21").Value2 = "496 (35.5)"
$ws.Range("C22").Value2 = "487 (34.8)"
$ws.Range("C23").Value2 = "182 (13.0)"
$ws.Range("C24").Value2 = "31 (2.2)"
$ws.Range("C25").Value2 = "83 (5.9)"
$ws.Range("C26").Value2 = "9 (0.6)"
$ws.Range("C27").Value2 = "51 (3.6)"

# MethylScoreAML Px Discovery cohort (column D) - 924 n
$ws.Range("D20").Value2 = "36 (3.9)"
$ws.Range("D21").Value2 = "491 (53.1)"
$ws.Range("D22").Value2 = "397 (43.0)"

# Validation cohort (column E) - 201 n
$ws.Range("E28").Value2 = "159 (79.1)"
$ws.Range("E29").Value2 = "42 (20.9)"
